$d = $word.ActiveDocument

# Locate the target paragraph ("Step 2 Konfigurasi Mongodb pada CMD") by
# its text content rather than a hard-coded index, so the script is
# resilient to the exact paragraph numbering.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text.TrimEnd()
    if ($t -eq "Step 2 Konfigurasi Mongodb pada CMD") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    Write-Output "ERROR: target paragraph 'Step 2 Konfigurasi Mongodb pada CMD' not found"
} else {
    $p = $d.Paragraphs.Item($targetIndex)
    $rng = $p.Range

    # Rebuild the whole paragraph via InsertXML so the exact run layout from
    # the target revision is produced: a new "(" run is added right after
    # "Step 2 ", the _GoBack bookmark moves to sit between "(" and
    # "Konfigurasi", and a new ")" run is appended after "pada CMD".
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p w:rsidR="00EC27B9" w:rsidRPr="00EC27B9" w:rsidRDefault="00EC27B9" w:rsidP="00EC27B9">' +
        '<w:pPr><w:ind w:firstLine="720"/></w:pPr>' +
        '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Step 2 </w:t></w:r>' +
        '<w:r><w:t>(</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t>Konfigurasi</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> Mongodb</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t>pada CMD</w:t></w:r>' +
        '<w:r><w:t>)</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xmlFrag)

    # InsertXML inserts the rebuilt content as a brand new paragraph placed
    # right before the original paragraph mark, which then survives as a
    # trailing empty paragraph. Remove that leftover paragraph mark so the
    # document keeps the same paragraph count/structure as before the edit.
    $insertedP = $d.Paragraphs.Item($targetIndex)
    $leftoverP = $d.Paragraphs.Item($targetIndex + 1)
    $delRng = $d.Range($insertedP.Range.End - 1, $leftoverP.Range.End)
    $delRng.Delete()

    Write-Output "Updated paragraph $targetIndex to: $($d.Paragraphs.Item($targetIndex).Range.Text)"
}
